$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fix the mistyped "test#kennect.io" values to "test@kennect.io"
$ws.Range("A3").Value = "test@kennect.io"
$ws.Range("A5").Value = "test@kennect.io"

# 2. Fix the mistyped "Qwerty@12345" values to "Qwerty@1234"
$ws.Range("B4").Value = "Qwerty@1234"
$ws.Range("B5").Value = "Qwerty@1234"

# 3. Remove all hyperlinks on the sheet, then re-create only the ones
#    that should remain: A2, B2 and A4.
$ws.Range("A1").Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:test@kennect.io") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:Qwerty@1234") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "mailto:test@kennect.io") | Out-Null

# Re-adding a hyperlink forces a fresh "Hyperlink" style onto the cell;
# put the cells back onto the original Hyperlink cell style.
$ws.Range("A2").Style = "Hyperlink"
$ws.Range("B2").Style = "Hyperlink"
$ws.Range("A4").Style = "Hyperlink"

# 4. Update the active selection shown when the workbook is opened.
$ws.Range("B5").Select() | Out-Null
